$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Median Value" (column C) and "Tier" (column D) figures, recalculated
# relative to the median AFTER merging with zip/census tract data.
$newData = @{
    2 = @{ C = 1.049645390070922; D = "4th Tier" }
    3 = @{ C = 0.9148936170212766; D = "Below Median" }
    4 = @{ C = 1.521276595744681; D = "1st Tier" }
    5 = @{ C = 1.035460992907801; D = "4th Tier" }
    6 = @{ C = 1.521276595744681; D = "1st Tier" }
    7 = @{ C = 1.893617021276596; D = "1st Tier" }
    8 = @{ C = 1.148936170212766; D = "3rd Tier" }
    9 = @{ C = 1.893617021276596; D = "1st Tier" }
    10 = @{ C = 0.1170212765957447; D = "Below Median" }
    11 = @{ C = 0.1170212765957447; D = "Below Median" }
    12 = @{ C = 0.1170212765957447; D = "Below Median" }
    13 = @{ C = 0.1170212765957447; D = "Below Median" }
    14 = @{ C = 0.1702127659574468; D = "Below Median" }
    15 = @{ C = 0.1702127659574468; D = "Below Median" }
    16 = @{ C = 0.425531914893617; D = "Below Median" }
    17 = @{ C = 0.5212765957446809; D = "Below Median" }
    18 = @{ C = 0.6382978723404256; D = "Below Median" }
    19 = @{ C = 0.4680851063829787; D = "Below Median" }
    20 = @{ C = 0.4680851063829787; D = "Below Median" }
    21 = @{ C = 0.5531914893617021; D = "Below Median" }
    22 = @{ C = 0.2553191489361702; D = "Below Median" }
    23 = @{ C = 0.2553191489361702; D = "Below Median" }
    24 = @{ C = 1.340425531914894; D = "2nd Tier" }
    25 = @{ C = 1.453900709219858; D = "2nd Tier" }
    26 = @{ C = 1.606382978723404; D = "1st Tier" }
    27 = @{ C = 1.606382978723404; D = "1st Tier" }
    28 = @{ C = 1.627659574468085; D = "1st Tier" }
    29 = @{ C = 1.134751773049645; D = "3rd Tier" }
    30 = @{ C = 1.390070921985816; D = "2nd Tier" }
    31 = @{ C = 0.925531914893617; D = "Below Median" }
    32 = @{ C = 1.340425531914894; D = "2nd Tier" }
    33 = @{ C = 1; D = "4th Tier" }
    34 = @{ C = 0.8404255319148937; D = "Below Median" }
    35 = @{ C = 0.7943262411347518; D = "Below Median" }
    36 = @{ C = 1.326241134751773; D = "3rd Tier" }
    37 = @{ C = 1.361702127659574; D = "2nd Tier" }
    38 = @{ C = 1.361702127659574; D = "2nd Tier" }
    39 = @{ C = 1.319148936170213; D = "3rd Tier" }
    40 = @{ C = 1.326241134751773; D = "2nd Tier" }
    41 = @{ C = 1.361702127659574; D = "2nd Tier" }
    42 = @{ C = 1.014184397163121; D = "4th Tier" }
    43 = @{ C = 1.340425531914894; D = "2nd Tier" }
    44 = @{ C = 1.056737588652482; D = "4th Tier" }
    45 = @{ C = 1.049645390070922; D = "4th Tier" }
    46 = @{ C = 0.3829787234042553; D = "Below Median" }
    47 = @{ C = 1; D = "4th Tier" }
    48 = @{ C = 0.8457446808510638; D = "Below Median" }
    49 = @{ C = 0.851063829787234; D = "Below Median" }
    50 = @{ C = 0.7943262411347518; D = "Below Median" }
    51 = @{ C = 1.049645390070922; D = "4th Tier" }
    52 = @{ C = 1; D = "4th Tier" }
    53 = @{ C = 1.308510638297872; D = "3rd Tier" }
    54 = @{ C = 0.4468085106382979; D = "Below Median" }
    55 = @{ C = 1.095744680851064; D = "3rd Tier" }
    56 = @{ C = 1.106382978723404; D = "3rd Tier" }
    57 = @{ C = 0.4609929078014184; D = "Below Median" }
    58 = @{ C = 0.3333333333333333; D = "Below Median" }
    59 = @{ C = 0.9787234042553191; D = "Below Median" }
    60 = @{ C = 1.24468085106383; D = "3rd Tier" }
    61 = @{ C = 1.148936170212766; D = "3rd Tier" }
    62 = @{ C = 1.319148936170213; D = "3rd Tier" }
    63 = @{ C = 1.014184397163121; D = "4th Tier" }
    64 = @{ C = 0.3333333333333333; D = "Below Median" }
    65 = @{ C = 1; D = "4th Tier" }
    66 = @{ C = 0.6879432624113476; D = "Below Median" }
    67 = @{ C = 0.4609929078014184; D = "Below Median" }
    68 = @{ C = 1.453900709219858; D = "2nd Tier" }
    69 = @{ C = 0.148936170212766; D = "Below Median" }
    70 = @{ C = 0.6808510638297872; D = "Below Median" }
    71 = @{ C = 1.453900709219858; D = "2nd Tier" }
    72 = @{ C = 0.1170212765957447; D = "Below Median" }
    73 = @{ C = 0.3191489361702128; D = "Below Median" }
    74 = @{ C = 0.8191489361702128; D = "Below Median" }
    75 = @{ C = 1.521276595744681; D = "1st Tier" }
    76 = @{ C = 1.453900709219858; D = "1st Tier" }
    77 = @{ C = 1.085106382978723; D = "3rd Tier" }
    78 = @{ C = 1.524822695035461; D = "1st Tier" }
    79 = @{ C = 0.9148936170212766; D = "Below Median" }
    80 = @{ C = 0.6666666666666666; D = "Below Median" }
}

foreach ($row in $newData.Keys) {
    $entry = $newData[$row]
    $ws.Cells.Item($row, 3).Value = $entry.C
    $ws.Cells.Item($row, 4).Value = $entry.D
}

Write-Output "Updated $($newData.Count) rows"
